# Auto-generated edit script: updates cryptos list values (prices/volumes) and
# re-orders a few coin rows (38-45) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $val)
    # Force the assigned value to be stored as literal text (matches the
    # source workbook where every data cell is an inline/shared string),
    # even when the text looks like a number (e.g. "1.000", "0.07704").
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.194.95"
Set-TextValue $ws.Range("E2") "  +3.56%  "
Set-TextValue $ws.Range("D3") "1.815.58"
Set-TextValue $ws.Range("E3") "  +5.04%  "
Set-TextValue $ws.Range("E4") "  -0.36%  "
Set-TextValue $ws.Range("D5") "329.95"
Set-TextValue $ws.Range("E5") "  +2.51%  "
Set-TextValue $ws.Range("D6") "1.000"
Set-TextValue $ws.Range("E6") "  -0.16%  "
Set-TextValue $ws.Range("D7") "0.4448"
Set-TextValue $ws.Range("E7") "  +5.80%  "
Set-TextValue $ws.Range("D8") "0.3698"
Set-TextValue $ws.Range("E8") "  +3.51%  "
Set-TextValue $ws.Range("D9") "44.61"
Set-TextValue $ws.Range("E9") "  -0.36%  "
Set-TextValue $ws.Range("D10") "0.07704"
Set-TextValue $ws.Range("E10") "  +3.85%  "
Set-TextValue $ws.Range("D11") "1.129"
Set-TextValue $ws.Range("E11") "  +2.01%  "
Set-TextValue $ws.Range("E12") "  -0.36%  "
Set-TextValue $ws.Range("D13") "22.15"
Set-TextValue $ws.Range("E13") "  +4.19%  "
Set-TextValue $ws.Range("D14") "6.293"
Set-TextValue $ws.Range("E14") "  +4.31%  "
Set-TextValue $ws.Range("D15") "7.586"
Set-TextValue $ws.Range("E15") "  +7.20%  "
Set-TextValue $ws.Range("D16") "1.834.50"
Set-TextValue $ws.Range("E16") "  +5.98%  "
Set-TextValue $ws.Range("D17") "92.90"
Set-TextValue $ws.Range("E17") "  +7.05%  "
Set-TextValue $ws.Range("D18") "0.00001085"
Set-TextValue $ws.Range("E18") "  +2.57%  "
Set-TextValue $ws.Range("D19") "0.06591"
Set-TextValue $ws.Range("E19") "  +10.56%  "
Set-TextValue $ws.Range("D20") "1.000"
Set-TextValue $ws.Range("E20") "  -0.16%  "
Set-TextValue $ws.Range("D21") "17.56"
Set-TextValue $ws.Range("E21") "  +5.30%  "
Set-TextValue $ws.Range("D22") "6.219"
Set-TextValue $ws.Range("E22") "  +3.00%  "
Set-TextValue $ws.Range("D23") "28.253.95"
Set-TextValue $ws.Range("E23") "  +3.51%  "
Set-TextValue $ws.Range("D24") "11.71"
Set-TextValue $ws.Range("E24") "  +3.89%  "
Set-TextValue $ws.Range("D25") "2.153"
Set-TextValue $ws.Range("E25") "  -10.09%  "
Set-TextValue $ws.Range("D26") "20.84"
Set-TextValue $ws.Range("E26") "  +4.28%  "
Set-TextValue $ws.Range("D27") "156.18"
Set-TextValue $ws.Range("E27") "  +4.58%  "
Set-TextValue $ws.Range("D28") "2.038.17"
Set-TextValue $ws.Range("E28") "  +5.75%  "
Set-TextValue $ws.Range("D29") "2.332"
Set-TextValue $ws.Range("E29") "  -0.24%  "
Set-TextValue $ws.Range("D30") "128.42"
Set-TextValue $ws.Range("E30") "  +2.17%  "
Set-TextValue $ws.Range("D31") "1.203"
Set-TextValue $ws.Range("E31") "  -2.42%  "
Set-TextValue $ws.Range("D32") "5.894"
Set-TextValue $ws.Range("E32") "  +6.01%  "
Set-TextValue $ws.Range("D33") "0.09230"
Set-TextValue $ws.Range("E33") "  +2.14%  "
Set-TextValue $ws.Range("D34") "3.655"
Set-TextValue $ws.Range("E34") "  -0.61%  "
Set-TextValue $ws.Range("D35") "13.12"
Set-TextValue $ws.Range("E35") "  +4.31%  "
Set-TextValue $ws.Range("D36") "0.02357"
Set-TextValue $ws.Range("E36") "  +3.98%  "
Set-TextValue $ws.Range("D37") "0.2183"
Set-TextValue $ws.Range("E37") "  +1.67%  "
Set-TextValue $ws.Range("B38") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D38") "5.185"
Set-TextValue $ws.Range("E38") "  +3.37%  "
Set-TextValue $ws.Range("B39") "Hedera"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D39") "0.06230"
Set-TextValue $ws.Range("E39") "  +2.10%  "
Set-TextValue $ws.Range("B40") "TheSandbox"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D40") "0.6581"
Set-TextValue $ws.Range("E40") "  +3.79%  "
Set-TextValue $ws.Range("B41") "FraxShare"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D41") "8.172"
Set-TextValue $ws.Range("E41") "  +3.82%  "
Set-TextValue $ws.Range("B42") "TrustWalletToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D42") "1.198"
Set-TextValue $ws.Range("E42") "  +1.51%  "
Set-TextValue $ws.Range("D43") "0.9995"
Set-TextValue $ws.Range("E43") "  -0.12%  "
Set-TextValue $ws.Range("B44") "WEMIXTOKEN"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D44") "1.401"
Set-TextValue $ws.Range("E44") "  -0.77%  "
Set-TextValue $ws.Range("B45") "EnergySwap"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D45") "13.95"
Set-TextValue $ws.Range("E45") "  +3.57%  "
Set-TextValue $ws.Range("D46") "0.6088"
Set-TextValue $ws.Range("E46") "  +5.08%  "
Set-TextValue $ws.Range("D47") "3.771"
Set-TextValue $ws.Range("E47") "  +1.16%  "
Set-TextValue $ws.Range("D48") "127.26"
Set-TextValue $ws.Range("E48") "  +2.59%  "
Set-TextValue $ws.Range("D49") "2.041"
Set-TextValue $ws.Range("E49") "  +5.90%  "
Set-TextValue $ws.Range("D50") "1.157"
Set-TextValue $ws.Range("E50") "  +6.13%  "
Set-TextValue $ws.Range("E51") "  +2.42%  "
